# 226091 Add create action for item, when product definition is updated
#
# Adds a new "create" action to the data-validation list on the "Items"
# sheet and appends a new item row (row 7) that uses it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# --- Append new row 7, mirroring the formatting/shape of row 3 -----------

$ws.Range("A7").Value = "ITM-1213-3316-0006"
$ws.Range("B7").Value = "Customer"
$ws.Range("C7").Value = "create"
$ws.Range("D7").Value = "65AB123BASD2"
$ws.Range("E7").Value = "NAV123456"
$ws.Range("F7").Value = "Description 6"
$ws.Range("G7").Value = "1m"
$ws.Range("H7").Value = "1y"
$ws.Range("I7").Value = "Published"
$ws.Range("J7").Value = "IGR-1213-3316-0002"
$ws.Range("K7").Value = "Default Group"
$ws.Range("L7").Value = "UNT-1916"
$ws.Range("M7").Value = "User"
$ws.Range("N7").Value = "False"
$ws.Range("O7").Value = "Migrate"
$ws.Range("P7").Value = "test ex 22"
$ws.Range("Q7").Value = 45292
$ws.Range("R7").Value = 45292

# Copy the number formatting (date format + font + alignment) from the
# equivalent cells in row 3 so row 7 re-uses the same style record
# instead of creating a new one.
$ws.Range("Q3:R3").Copy()
$ws.Range("Q7:R7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Extend the dropdown list used for the Action column -----------------

$dv = $ws.Range("C2:C1048576").Validation
$dv.Formula1 = """-,create,update,review,publish,unpublish"""

# --- Update the active selection, matching the sheet state after editing -

[void]$ws.Activate()
[void]$ws.Range("C13").Select()
